$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handback" -- refresh the report rows so the two
# e2e markdown files (and their derived xliff / timestamp columns) point at
# the newly generated handback artifacts instead of the old ones.
# ---------------------------------------------------------------------------

# New cell values, keyed by sheet name / cell address.
$cellUpdates = @{
    "Overview" = @{
        "A2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "B2" = "e2e\718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "G2" = "2016-08-24 03:02:15"
        "A3" = "fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
        "B3" = "e2e\fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
        "G3" = "2016-08-24 03:02:15"
    }
    "zh-cn" = @{
        "A2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "G2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.5639ba5301a44a582f0cdbc9075255218643f13d.zh-cn.xlf"
        "H2" = "2016-08-24 03:02:10"
        "I2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "J2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.5639ba5301a44a582f0cdbc9075255218643f13d.zh-cn.xlf"
        "K2" = "2016-08-24 03:02:35"
        "A3" = "fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
        "G3" = "718087ba-2aa5-46e7-a393-b22be688e1b5.5639ba5301a44a582f0cdbc9075255218643f13d.zh-cn.xlf"
        "H3" = "2016-08-24 03:02:10"
        "I3" = "fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
        "J3" = "718087ba-2aa5-46e7-a393-b22be688e1b5.5639ba5301a44a582f0cdbc9075255218643f13d.zh-cn.xlf"
        "K3" = "2016-08-24 03:02:35"
    }
    "de-de" = @{
        "A2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "G2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.5639ba5301a44a582f0cdbc9075255218643f13d.de-de.xlf"
        "H2" = "2016-08-24 03:02:15"
        "I2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "J2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.5639ba5301a44a582f0cdbc9075255218643f13d.de-de.xlf"
        "K2" = "2016-08-24 03:02:43"
        "A3" = "fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
        "G3" = "718087ba-2aa5-46e7-a393-b22be688e1b5.5639ba5301a44a582f0cdbc9075255218643f13d.de-de.xlf"
        "H3" = "2016-08-24 03:02:15"
        "I3" = "fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
        "J3" = "718087ba-2aa5-46e7-a393-b22be688e1b5.5639ba5301a44a582f0cdbc9075255218643f13d.de-de.xlf"
        "K3" = "2016-08-24 03:02:43"
    }
}

# New hyperlink display text, keyed by sheet name / cell address (the
# underlying hyperlink target URLs themselves are unchanged).
$hyperlinkUpdates = @{
    "Overview" = @{
        "B2" = "e2e\718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "B3" = "e2e\fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
    }
    "zh-cn" = @{
        "A2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "I2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "A3" = "fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
        "I3" = "fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
    }
    "de-de" = @{
        "A2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "I2" = "718087ba-2aa5-46e7-a393-b22be688e1b5.md"
        "A3" = "fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
        "I3" = "fffff26fed17-6f06-4a6b-9c3b-042bfa3759af.md"
    }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name

    if ($cellUpdates.ContainsKey($name)) {
        $updates = $cellUpdates[$name]
        foreach ($addr in $updates.Keys) {
            $ws.Range($addr).Value = $updates[$addr]
        }
    }

    if ($hyperlinkUpdates.ContainsKey($name)) {
        $linkUpdates = $hyperlinkUpdates[$name]
        foreach ($hl in $ws.Hyperlinks) {
            $addr = $hl.Range.Address()
            foreach ($cellAddr in $linkUpdates.Keys) {
                $fullAddr = '$' + $cellAddr.Substring(0,1) + '$' + $cellAddr.Substring(1)
                if ($addr -eq $fullAddr) {
                    $hl.TextToDisplay = $linkUpdates[$cellAddr]
                }
            }
        }
    }
}
